$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E ("ABRIL"-style): D/56 + D -------------------------------
$ws.Range("E3").Formula = "=D3/56+D3"
$ws.Range("E4:E12").Formula = "=D4/56+D4"
$ws.Range("E14").Formula = "=AVERAGE(E3:E12)"
$ws.Range("E17").FormulaArray = "=SUM(G14+E3:E12)"

# --- New column F ("MAYO"-style): SUM(D,D,D) ------------------------------
$ws.Range("F3").Formula = "=SUM(D3,D3,D3)"
$ws.Range("F4:F12").Formula = "=SUM(D4,D4,D4)"
$ws.Range("F14").Formula = "=AVERAGE(F3:F12)"
$ws.Range("F17").FormulaArray = "=SUM(H14+F3:F12)"

# --- Formatting: column E gets the built-in Percent cell style -----------
# (applied after all formulas are written so the style doesn't get
#  auto-propagated onto cells that merely reference this range)
$ws.Range("E3:E12").Style = "Percent"

# --- Selection moved from E3 to G3 ----------------------------------------
$ws.Range("G3").Select() | Out-Null
